$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 3 (R row) with Week 17 data ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 204
$wsOff.Range("C3").Value = 142
$wsOff.Range("D3").Value = 58
$wsOff.Range("E3").Value = 26

# --- DEF sheet: update row 3 (R row) with Week 17 data ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 267
$wsDef.Range("C3").Value = 208
$wsDef.Range("D3").Value = 41
$wsDef.Range("E3").Value = 19
$wsDef.Range("G3").Value = 3
